# Update the answers in the "two-digit number divided by one-digit number"
# worksheet table. Each data row of the table (rows 1, 5, 9, 13, 17) holds
# five answer strings across its five columns; replace each cell's text
# with the corresponding new value while preserving the existing run
# formatting (font, size, etc.).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, column) -> new text, in the exact order/values from the diff.
$updates = @(
    @{ Row = 1;  Col = 1; New = "68÷7=9, 5" }
    @{ Row = 1;  Col = 2; New = "74÷3=24, 2" }
    @{ Row = 1;  Col = 3; New = "44÷7=6, 2" }
    @{ Row = 1;  Col = 4; New = "46÷3=15, 1" }
    @{ Row = 1;  Col = 5; New = "15÷4=3, 3" }

    @{ Row = 5;  Col = 1; New = "54÷9=6, 0" }
    @{ Row = 5;  Col = 2; New = "53÷3=17, 2" }
    @{ Row = 5;  Col = 3; New = "22÷2=11, 0" }
    @{ Row = 5;  Col = 4; New = "42÷3=14, 0" }
    @{ Row = 5;  Col = 5; New = "21÷7=3, 0" }

    @{ Row = 9;  Col = 1; New = "85÷9=9, 4" }
    @{ Row = 9;  Col = 2; New = "63÷6=10, 3" }
    @{ Row = 9;  Col = 3; New = "17÷7=2, 3" }
    @{ Row = 9;  Col = 4; New = "34÷7=4, 6" }
    @{ Row = 9;  Col = 5; New = "11÷7=1, 4" }

    @{ Row = 13; Col = 1; New = "33÷2=16, 1" }
    @{ Row = 13; Col = 2; New = "38÷6=6, 2" }
    @{ Row = 13; Col = 3; New = "18÷8=2, 2" }
    @{ Row = 13; Col = 4; New = "80÷5=16, 0" }
    @{ Row = 13; Col = 5; New = "18÷9=2, 0" }

    @{ Row = 17; Col = 1; New = "15÷6=2, 3" }
    @{ Row = 17; Col = 2; New = "35÷8=4, 3" }
    @{ Row = 17; Col = 3; New = "57÷4=14, 1" }
    @{ Row = 17; Col = 4; New = "95÷4=23, 3" }
    @{ Row = 17; Col = 5; New = "38÷6=6, 2" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Exclude the trailing cell-mark / paragraph-mark character so only the
    # visible text is replaced, keeping the run's formatting intact.
    $rng.End = $rng.End - 1
    $rng.Text = $u.New
}

Write-Host "Updated $($updates.Count) cells"
